# Add the "only_zero_score" column.
#
# The header label "only_zero_score" is written to S1, and a per-row
# yes/no result (whether every score column for all rows sharing the same
# input_id is zero/blank) is written into column T for each data row
# (T2:T23) - matching the layout produced by the original author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("S1").Value = "only_zero_score"

# Score columns used to determine whether a person's candidate rows are
# "only zero score": match_score(D), altname_score(E), biogaddr_score(G),
# entry_score(I), posting_score(L), kin_score(N), death_nh_score(P).
$scoreCols = @("D", "E", "G", "I", "L", "N", "P")

$lastRow = 23
$inputIdCol = "A"

for ($row = 2; $row -le $lastRow; $row++) {
    $inputId = $ws.Range("$inputIdCol$row").Value2

    $anyNonZero = $false
    for ($r2 = 2; $r2 -le $lastRow; $r2++) {
        if ($ws.Range("$inputIdCol$r2").Value2 -eq $inputId) {
            foreach ($col in $scoreCols) {
                $val = $ws.Range("$col$r2").Value2
                if ($val -ne $null -and $val -ne 0) {
                    $anyNonZero = $true
                }
            }
        }
    }

    if ($anyNonZero) {
        $ws.Range("T$row").Value = "no"
    } else {
        $ws.Range("T$row").Value = "yes"
    }
}

Write-Output "only_zero_score column created"
